# spp_gp_representativeness.xlsx - add "fish" representativeness rows
# based on Mike Kingsford input, per commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Window / view cosmetics (best effort - some properties may be
#    no-ops on this host, which is fine, they are purely cosmetic).
# ---------------------------------------------------------------------
try { $excel.ActiveWindow.Height = 18000 } catch {}

$ws.Activate()

# ---------------------------------------------------------------------
# 2) New data rows 155:176 - fish representativeness entries.
#    Columns: A = taxon ("fish"), B = species, C = representative
#    rank ("genus"), D = number of species in the genus.
# ---------------------------------------------------------------------

$species = @(
    @{Row=155; Name="Acanthopagrus australis";     N=9;  Style="ValignFillA"},
    @{Row=156; Name="Amphiprion akindynos";         N=28; Style="Valign"},
    @{Row=157; Name="Brachionichthys hirsutus";     N=5;  Style="Valign"},
    @{Row=158; Name="Cephalopholis cyanostigmata";  N=9;  Style="Valign"},
    @{Row=159; Name="Clupea harengus";              N=4;  Style="Valign"},
    @{Row=160; Name="Corphaena hippurus";           N=2;  Style="Valign"},
    @{Row=161; Name="Chromis atripectoralis";       N=73; Style="ValignFillA"},
    @{Row=162; Name="Chromis dispilus";             N=3;  Style="ThemeFillA"},
    @{Row=163; Name="Engraulis ringens";            N=50; Style="Valign"},
    @{Row=164; Name="Epinephelus coioides";         N=6;  Style="ThemeSmallFillA"},
    @{Row=165; Name="Epinephelus fasciatus";        N=11; Style="ValignFillA"},
    @{Row=166; Name="Hoplostethus japonicus";       N=13; Style="Theme"},
    @{Row=167; Name="Molva dypterygia";             N=7;  Style="ValignFillB"},
    @{Row=168; Name="Oncorhynchus nerka";           N=18; Style="Valign"},
    @{Row=169; Name="Pagrus auratus";               N=10; Style="ValignFillA"},
    @{Row=170; Name="Pomacentrus amboinensis";      N=53; Style="Valign"},
    @{Row=171; Name="Parma microlepis";             N=10; Style="Valign"},
    @{Row=172; Name="Plectropomus leopardus";       N=2;  Style="Valign"},
    @{Row=173; Name="Plectropomus areolatus";       N=2;  Style="Valign"},
    @{Row=174; Name="Rexea solandri";               N=6;  Style="Valign"},
    @{Row=175; Name="Spratelloides delicatulus";    N=9;  Style="Valign"},
    @{Row=176; Name="Thunnus maccoyii";             N=8;  Style="Valign"}
)

$fillA = 3407769   # RGB(153,255,51) -> FF99FF33 (light green)
$fillB = 6750054   # RGB(102,255,102) -> FF66FF66 (darker green)

foreach ($s in $species) {
    $r = $s.Row

    $rngA = $ws.Range("A$r")
    $rngB = $ws.Range("B$r")
    $rngC = $ws.Range("C$r")
    $rngD = $ws.Range("D$r")

    $rngA.Value = "fish"
    $rngB.Value = $s.Name
    $rngC.Value = "genus"
    $rngD.Value = $s.N

    # All rows: species name (B) is italic.
    $rngB.Font.Italic = $true

    # All rows: representative-rank (C) cell has the light-green fill.
    $rngC.Interior.Color = $fillA

    # All rows: number-of-species (D) cell is horizontally centered.
    $rngD.HorizontalAlignment = -4108   # xlCenter

    switch ($s.Style) {
        "Valign" {
            $rngB.Font.Color = 0
            $rngB.VerticalAlignment = -4108   # xlCenter
        }
        "ValignFillA" {
            $rngB.Font.Color = 0
            $rngB.VerticalAlignment = -4108
            $rngB.Interior.Color = $fillA
        }
        "ValignFillB" {
            $rngB.Font.Color = 0
            $rngB.VerticalAlignment = -4108
            $rngB.Interior.Color = $fillB
        }
        "ThemeFillA" {
            $rngB.Font.ThemeColor = 1
            $rngB.Interior.Color = $fillA
            $rngD.Font.Italic = $true
            $rngD.Font.ThemeColor = 1
        }
        "ThemeSmallFillA" {
            $rngB.Font.ThemeColor = 1
            $rngB.Font.Size = 10
            $rngB.Interior.Color = $fillA
        }
        "Theme" {
            $rngB.Font.ThemeColor = 1
        }
    }
}

# ---------------------------------------------------------------------
# 3) Sheet view - scroll position, zoom and final selection, matching
#    where the author ended up after entering the new fish block.
# ---------------------------------------------------------------------
try { $excel.ActiveWindow.ScrollRow = 148 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
try { $excel.ActiveWindow.Zoom = 150 } catch {}

$ws.Range("C155:C176").Select()

Write-Host "fish representativeness rows added"
